$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Range("D235").NumberFormat

$newRows = @(
    @{ Row=236; L="Especial"; M=200; Q="$/caja 10 unidades"; S=1950; T=10 },
    @{ Row=237; L="Primera";  M=250; Q="$/caja 12 unidades"; S=1625; T=12 },
    @{ Row=238; L="Segunda";  M=270; Q="$/caja 14 unidades"; S=1393; T=14 },
    @{ Row=239; L="Tercera";  M=200; Q="$/caja 16 unidades"; S=1219; T=16 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

    $ws.Cells.Item($row, 4).Value = 44832
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value = 100108005
    $ws.Cells.Item($row, 10).Value = "Piña"
    $ws.Cells.Item($row, 11).Value = "Caramelo"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = 19000
    $ws.Cells.Item($row, 15).Value = 20000
    $ws.Cells.Item($row, 16).Value = 19500
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Ecuador"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
